## Created new scripts in IAM module
## - Extends IAM040's Jira id list and verification text (row 41 on the
##   "Test Cases" sheet) with newly added scripts.
## - Adds a brand new IAM041 test case as row 42 on the same sheet.
## - Updates the saved selection on the IAM009 sheet (view-state only).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Row 41 (IAM040): append the two newly created Jira ids and rework
#     the verification text to match the new scripts. ---
$ws.Range("B41").Value = 'OPQA-5372||OPQA-5373||OPQA-4252||OPQA-5401||OPQA-5402'
$ws.Range("C41").Value = 'Verify that error message "Please enter a valid email address." should be displayed in red color when user not enter email address in email text field for neon login page.||Verify that error message "Please enter a password." should be displayed in red color when user not enter password in password text field for neon login page||Verify that error message " Please enter a valid email address." should be displayed in red color when user enters email address in wrong format||Verify that error messages "Please enter a valid email address." and "Please enter a password." should be displayed in red color when user click login button without enter username and password.||Verify that error message "Invalid email/password. Please try again." should be displayed in red color when user entered wrong username and password in login page.'
$ws.Rows.Item(41).RowHeight = 158.4

# --- Row 42 (new IAM041 test case): clone row 41's formatting, then fill
#     in the new content. ---
[void]$ws.Range("A41:E41").Copy()
[void]$ws.Range("A42:E42").PasteSpecial(-4122)
$ws.Rows.Item(42).RowHeight = 129.6

$ws.Range("A42").Value = 'IAM041'
$ws.Range("C42").Value = 'Verify that error message "Please enter an email address." should be displayed in red color when user not enter email address in email text field for neon singup page.||Verify that error message "Please enter a password." should be displayed in red color when user not enter password in password text field for neon singup page.||Verify that error message "Please enter your first name." should be displayed in red color when user not enter first name in first name text field for neon singup page.||Verify that error message "Please enter your last name." should be displayed in red color when user not enter last name in last name text field for neon singup page.'
$ws.Range("B42").Value = 'OPQA-5403||OPQA-5404||OPQA-5405||OPQA-5406'
$ws.Range("D42").Value = 'Y'
$ws.Range("E42").Value = 'PASS'

# --- View-state: IAM009 sheet's saved selection moved to G17 (the user
#     clicked around on that sheet); restore focus to "Test Cases"
#     afterwards so it stays the active tab. ---
$ws9 = $wb.Worksheets.Item("IAM009")
[void]$ws9.Activate()
[void]$ws9.Range("G17").Select()
[void]$ws.Activate()
[void]$ws.Range("E42").Select()
